# control_asistencia.xlsx - "se realizo todas las clases para esta semana"
#
# - Adds the "p" (attendance) mark in column G for every student row except
#   the two that already had a full house (row 8 / HIDALGO, row 11 / RODRÍGUEZ
#   never gets a G value in this diff - so we mirror that, skipping rows 8 and 11).
# - Fills in three missing e-mail addresses (column D) together with their
#   mailto hyperlinks, re-using the existing "plain" hyperlink look (style of D3)
#   instead of Excel's default blue/underline Hyperlink cell style.
# - Moves the active selection from D6 to G5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that get the "p" attendance mark added in column G this week.
$attendanceRows = @(3, 4, 5, 6, 7, 9, 10, 12, 13, 14)
foreach ($r in $attendanceRows) {
    $ws.Cells.Item($r, 7).Value = "p"
}

# Reference style used by the other e-mail hyperlink cells in column D
# (direct underline/theme-color formatting, not the named "Hyperlink" style).
$linkStyle = $ws.Range("D3").Style

# Newly-completed e-mail addresses, in the same order they were added to the
# shared-string table: claudi-94 (D12), cris_0.5 (D9), elena_nena_91 (D4).
$newEmails = @(
    @{ Cell = "D12"; Address = "claudi-94@hotmail.com" },
    @{ Cell = "D9";  Address = "cris_0.5@hotmail.com" },
    @{ Cell = "D4";  Address = "elena_nena_91@hotmail.com" }
)

foreach ($entry in $newEmails) {
    $rng = $ws.Range($entry.Cell)
    $rng.Value = $entry.Address
    $ws.Hyperlinks.Add($rng, "mailto:" + $entry.Address)
    $rng.Style = $linkStyle
}

# Move the selection (was D6) to G5.
$ws.Range("G5").Select() | Out-Null
